# Applies the "Implemented targeting for enemies" commit:
#  - Cleans up grammar-checker proofErr markers / run-splits on several
#    bullet points whose visible text does not actually change.
#  - Rewrites the "Make find target..." bullet to
#    "Implement enemy squad forming/merging".
#  - Fills in the previously-empty trailing bullet with
#    "Fix enemies targeting followers before the home fire" (kept as two
#    separate runs, matching the source edit) and appends a new bullet
#    "Follower retaliation (target enemies  if hit or if squad member is hit)".

$d = $word.ActiveDocument

# --- 1) "Spawners appear in corrupted lands, ..." -------------------------
# Text is unchanged; only the gramStart/gramEnd proofErr-split runs around
# "lands," collapse back into a single run. A same-text Find/Replace merges
# the runs and drops the proofErr markers that sit fully inside the match.
[void]$d.Content.Find.Execute(
    "Spawners appear in corrupted lands, the player must find them and destroy them before they grow",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Spawners appear in corrupted lands, the player must find them and destroy them before they grow",
    2)

# --- 2) "More enemies and tougher enemies spawn as the player progresses" -
# Here the trailing proofErr (gramEnd) sits right at the paragraph mark. A
# same-text Find/Replace leaves that dangling gramEnd behind (it has no
# character width, so it isn't "inside" the matched text). Instead, replace
# this paragraph together with the one right before it (reproduced as-is)
# via a raw OOXML insertion spanning both paragraphs' marks, which cleanly
# drops every proofErr in the replaced span.
$paras = $d.Paragraphs
$prevPara = $paras.Item(21)
$thisPara = $paras.Item(22)
if ($prevPara.Range.Text.TrimEnd([char]13) -ne "Enemies target the home building (firepit etc.) until they are hit or are blocked by a wall") {
    throw "unexpected neighbour paragraph text; aborting to avoid corrupting the document"
}
$spanRange = $d.Range($prevPara.Range.Start, $thisPara.Range.End)
$twoParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t>Enemies target the home building (firepit etc.) until they are hit or are blocked by a wall</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t>More enemies and tougher enemies spawn as the player progresses</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$spanRange.InsertXML($twoParaXml)

# --- 3) "Squads - allows soldiers, archers and priests ..." ---------------
[void]$d.Content.Find.Execute(
    "Squads " + [char]0x2013 + " allows soldiers, archers and priests to be commanded as groups rather than as single units. Combine squads to create armies",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Squads " + [char]0x2013 + " allows soldiers, archers and priests to be commanded as groups rather than as single units. Combine squads to create armies",
    2)

# --- 4) "Explosive - moves really slowly but deals AOE damage ..." --------
[void]$d.Content.Find.Execute(
    "Explosive " + [char]0x2013 + " moves really slowly but deals AOE damage when in range of followers or buildings. They must be killed in time",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Explosive " + [char]0x2013 + " moves really slowly but deals AOE damage when in range of followers or buildings. They must be killed in time",
    2)

# --- 5) "Make find target ..." -> "Implement enemy squad forming/merging" -
[void]$d.Content.Find.Execute(
    "Make find target first check if the previous target was in a squad " + [char]0x2013 + " if so find the closest enemy in the squad",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implement enemy squad forming/merging",
    2)

# --- 6) Fill the trailing empty bullet + append a new bullet --------------
# The last paragraph in the document is an empty ListParagraph/numId=14
# bullet. Replace it (and add a new bullet after it) via a raw OOXML
# insertion so the two sentences of bullet "Fix enemies ..." stay as two
# distinct runs, exactly like the source edit.
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$insertionPoint = $lastPara.Range

$openXmlPackage = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr>
<w:r><w:t xml:space="preserve">Fix enemies </w:t></w:r>
<w:r><w:t>targeting followers before the home fire</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr>
<w:r><w:t>Follower retaliation (target enemies  if hit or if squad member is hit)</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$insertionPoint.InsertXML($openXmlPackage)

# InsertXML adds the two new bullets *before* the (still empty) original
# bullet paragraph, leaving a redundant empty trailing paragraph mark.
# Remove that redundant paragraph mark by deleting the range spanning from
# just before it (end of the newly added "Follower retaliation" bullet)
# through its own end, merging it away.
$paras2 = $d.Paragraphs
$trailingEmpty = $paras2.Item($paras2.Count)
$newSecondBullet = $paras2.Item($paras2.Count - 1)
$mergeRange = $d.Range($newSecondBullet.Range.End - 1, $trailingEmpty.Range.End)
[void]$mergeRange.Delete()
